$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 157, pushing the
# existing data (rows 157-175) down to rows 159-177.
$ws.Rows(157).Insert()
$ws.Rows(157).Insert()

# --- New row 157 ---
$ws.Range("A157").Value = 1
$ws.Range("B157").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C157").Value = "Arica y Parinacota"
$ws.Range("D157").Value = 44918
$ws.Range("E157").Value = 15
$ws.Range("F157").Value = "Fruta"
$ws.Range("G157").Value = 100108
$ws.Range("H157").Value = "Tropicales y subtropicales"
$ws.Range("I157").Value = 100108002
$ws.Range("J157").Value = "Mango"
$ws.Range("K157").Value = "Sin especificar"
$ws.Range("L157").Value = "Especial"
$ws.Range("M157").Value = 750
$ws.Range("N157").Value = 4500
$ws.Range("O157").Value = 5000
$ws.Range("P157").Value = 4833
$ws.Range("Q157").Value = "$/bandeja 4 kilos"
$ws.Range("R157").Value = "Perú"
$ws.Range("S157").Value = 1208
$ws.Range("T157").Value = 4

# --- New row 158 ---
$ws.Range("A158").Value = 1
$ws.Range("B158").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C158").Value = "Arica y Parinacota"
$ws.Range("D158").Value = 44918
$ws.Range("E158").Value = 15
$ws.Range("F158").Value = "Fruta"
$ws.Range("G158").Value = 100108
$ws.Range("H158").Value = "Tropicales y subtropicales"
$ws.Range("I158").Value = 100108002
$ws.Range("J158").Value = "Mango"
$ws.Range("K158").Value = "Sin especificar"
$ws.Range("L158").Value = "Primera"
$ws.Range("M158").Value = 500
$ws.Range("N158").Value = 4500
$ws.Range("O158").Value = 5000
$ws.Range("P158").Value = 4750
$ws.Range("Q158").Value = "$/bandeja 4 kilos"
$ws.Range("R158").Value = "Perú"
$ws.Range("S158").Value = 1188
$ws.Range("T158").Value = 4
